# Commit: "uniformised load date/ header date and filename /file between
# globalscape and swift (kept swift format)"
#
# The GLOB_subs sheet used its own GlobalScape-specific field names
# ("loaddate" / "fileid") at the bottom of the Variable/Type list. This
# change renames them to match the SWIFT naming already used elsewhere in
# the workbook ("header_date" / "file_id"), and adds the matching Type
# value ("character") for file_id/fileid, which GLOB_subs was previously
# missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLOB_subs")
[void]$ws.Activate()

# loaddate -> header_date (Type column, B48, stays "date")
$ws.Range("A48").Value = "header_date"

# fileid -> file_id, and give it the "character" type (B49), matching the
# SWIFT_subs sheet's file_id/character pairing
$ws.Range("A49").Value = "file_id"
$ws.Range("B49").Value = "character"

# Selection left on the last-edited cell
[void]$ws.Range("A49").Select()
